$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing value: PopSizeDE (G2) 10 -> 30 ---
$ws.Range("G2").Value = 30

# --- Remove the fill style from K1/L1 header cells (keep the border) ---
$ws.Range("K1:L1").Interior.Pattern = -4142
$ws.Range("K1:L1").Borders.Item(9).LineStyle = 1

# --- Add the new "objective selection" columns ---
# Write cells in an order that reproduces the target sharedStrings sequence:
# OF_options, TAC, GHG, CAP, COP, OFs
$ws.Range("O1").Value = "OF_options"
$ws.Range("N2").Value = "TAC"
$ws.Range("O2").Value = "TAC"
$ws.Range("N3").Value = "GHG"
$ws.Range("P3").Value = "GHG"
$ws.Range("O3").Value = "CAP"
$ws.Range("P2").Value = "COP"
$ws.Range("N1").Value = "OFs"

# --- Restore the selection/active cell as in the edited workbook ---
$ws.Range("Q8").Select()
